$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 0.6265618270263076
$ws.Range("H4").Value = 1.596011943469663
$ws.Range("I4").Value = 0.002366712753428146
$ws.Range("J4").Value = 0.0005064628257969161
$ws.Range("L4").Value = 0.0006180772033985704
$ws.Range("M4").Value = 0.0009185673316096654
